# Applies the row re-shuffle / data update described in the commit diff
# to the "domain_comparison_filter" worksheet.
#
# The underlying OOXML diff shows the shared-strings table being reordered
# (moving the Ig/COG1100/CUB/Ion_trans blocks earlier, dropping the SPEC
# entry, and adding a new zf-H2C2_2 entry at the end) together with a
# handful of explicit cell edits on top of that reorder. The net visible
# effect on the worksheet grid is the set of cell updates below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 'Ig'
$ws.Range("B6").Value = 'common'
$ws.Range("C6").Value = 'Immunoglobulin domain'

$ws.Range("A7").Value = 'COG1100'
$ws.Range("B7").Value = 'KRAS'
$ws.Range("C7").Value = 'GTPase SAR1 and related small G proteins [General function prediction only]'

$ws.Range("A8").Value = 'CUB'
$ws.Range("B8").Value = 'common'
$ws.Range("C8").Value = 'CUB domain'

$ws.Range("A9").Value = 'I-set'
$ws.Range("B9").Value = 'common'
$ws.Range("C9").Value = 'Immunoglobulin I-set domain'

$ws.Range("A10").Value = 'LamG'
$ws.Range("B10").Value = 'common'
$ws.Range("C10").Value = 'Laminin G domain'

$ws.Range("A11").Value = 'Ion_trans'
$ws.Range("B11").Value = 'common'
$ws.Range("C11").Value = 'Ion transport protein'

$ws.Range("A12").Value = 'COG2319'
$ws.Range("B12").Value = 'common'
$ws.Range("C12").Value = 'FOG: WD40 repeat [General function prediction only]'

$ws.Range("A14").Value = 'Tryp_SPc'
$ws.Range("B14").Value = 'common'
$ws.Range("C14").Value = 'Trypsin-like serine protease'

$ws.Range("A15").Value = '7tm_4'
$ws.Range("B15").Value = 'common'
$ws.Range("C15").Value = 'Olfactory receptor'

$ws.Range("A16").Value = 'ANK'
$ws.Range("B16").Value = 'common'
$ws.Range("C16").Value = 'ankyrin repeats'

$ws.Range("A17").Value = 'S_TKc'
$ws.Range("B17").Value = 'KRAS'
$ws.Range("C17").Value = 'Serine/Threonine protein kinases, catalytic domain'

$ws.Range("A18").Value = 'LIC'
$ws.Range("B18").Value = 'KRAS'
$ws.Range("C18").Value = 'Cation transporter family protein'

$ws.Range("A20").Value = 'PTPc'
$ws.Range("B20").Value = 'KRAS'
$ws.Range("C20").Value = 'Protein tyrosine phosphatases (PTP) catalyze the dephosphorylation of phosphotyrosine peptides'

$ws.Range("A21").Value = 'MYSc'
$ws.Range("B21").Value = 'common'
$ws.Range("C21").Value = 'Myosin. Large ATPases'

$ws.Range("A24").Value = 'P53'
$ws.Range("B24").Value = 'EGFR'
$ws.Range("C24").Value = 'P53 DNA-binding domain'

$ws.Range("A25").Value = 'PDZ_signaling'
$ws.Range("B25").Value = 'EGFR'
$ws.Range("C25").Value = 'PDZ domain found in a variety of Eumetazoan signaling molecules, often in tandem arrangements. May be responsible for specific protein-protein interactions, as most PDZ domains bind C-terminal polypeptides, and binding to internal (non-C-terminal)...'

$ws.Range("A26").Value = 'zf-H2C2_2'
$ws.Range("B26").Value = 'EGFR'
$ws.Range("C26").Value = 'Zinc-finger double domain'
